# Apply "Add data for 2022-06-26" update to the carjacking-by-neighborhood-by-month workbook.
# This bumps the "through" date in the sheet name / header label by one day,
# and increments/adds the relevant monthly totals that changed as a result
# of the newly-added day of raw data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the column header label that both
# reference the "through" date.
$ws.Name = "Through 2022-06-18"
$ws.Range("B1").Value = "June 2022 (through June 18)"

# Cell value updates: Address -> New Value
$updates = @{
    "AR3"  = 1
    "H4"   = 4
    "B5"   = 1
    "AF5"  = 5
    "Z6"   = 2
    "N9"   = 3
    "T10"  = 2
    "B12"  = 5
    "AF13" = 2
    "H14"  = 8
    "N14"  = 5
    "B17"  = 1
    "AF30" = 2
    "B33"  = 1
    "B39"  = 1
    "B48"  = 1
    "N51"  = 2
    "B68"  = 3
    "H68"  = 4
    "B70"  = 4
    "B85"  = 1
    "B89"  = 2
    "AF96" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
